$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(1, 2).Value = "year"
$ws.Cells.Item(1, 3).Value = "description"
$ws.Cells.Item(1, 4).Value = "link"
$ws.Cells.Item(1, 5).Value = "img"
$ws.Cells.Item(1, 6).Value = "image_link"
$ws.Cells.Item(1, 7).Value = "branch"
$ws.Cells.Item(1, 8).Value = "city"
$ws.Cells.Item(1, 9).Value = "state"
$ws.Cells.Item(1, 10).Value = "about"

$ws.Cells.Item(2, 1).Value = "Aditya Rana"
$ws.Cells.Item(2, 2).Value = "first"
$ws.Cells.Item(2, 3).Value = "Volunteer Member"
$ws.Cells.Item(2, 4).Value = "https://www.linkedin.com/in/aditya-rana-6156071aa"
$ws.Cells.Item(2, 5).Value = "../members/aditiya.webp"
$ws.Cells.Item(2, 6).Value = "1LjuDNVhaQmaOaujbP3XWpNtAaXtnAZNo"
$ws.Cells.Item(2, 7).Value = "Computer science and Engineering"
$ws.Cells.Item(2, 8).Value = "Kangra"
$ws.Cells.Item(2, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(2, 10).Value = "I am an instance of my own class | Fresher @ CSE NITH | Always learning and implementing latest and greatest technologies | Freelancer @ fiverr | Android Developer |"

$ws.Cells.Item(3, 1).Value = "Armaan Shukla"
$ws.Cells.Item(3, 2).Value = "first"
$ws.Cells.Item(3, 3).Value = "Volunteer Member"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "../members/armaan.webp"
$ws.Cells.Item(3, 6).Value = "1Tpnk2k0ZIv73FWy7JeLxcRBjC_Oq7O-f"
$ws.Cells.Item(3, 7).Value = "Mathematics and Computing"
$ws.Cells.Item(3, 8).Value = "Kangra"
$ws.Cells.Item(3, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(3, 10).Value = "Batch 2025"

$ws.Cells.Item(4, 1).Value = "Chandan Kumar"
$ws.Cells.Item(4, 2).Value = "first"
$ws.Cells.Item(4, 3).Value = "Volunteer Member"
$ws.Cells.Item(4, 4).Value = "https://www.linkedin.com/in/chandan-kumar-19880a22a/"
$ws.Cells.Item(4, 5).Value = "../members/chandan.webp"
$ws.Cells.Item(4, 6).Value = "1522oz7z6THBgb0FxF4DmHlyI3YQpviNM"
$ws.Cells.Item(4, 7).Value = "Electronics and Communication Engineering"
$ws.Cells.Item(4, 8).Value = "Varanasi"
$ws.Cells.Item(4, 9).Value = "Uttar Pradesh"
$ws.Cells.Item(4, 10).Value = "Batch 2025"

$ws.Cells.Item(5, 1).Value = "Charu"
$ws.Cells.Item(5, 2).Value = "first"
$ws.Cells.Item(5, 3).Value = "Volunteer Member"
$ws.Cells.Item(5, 4).Value = "https://www.linkedin.com/in/charu-229665223"
$ws.Cells.Item(5, 5).Value = "../members/charu.webp"
$ws.Cells.Item(5, 6).Value = "1gYhb-XoUU6OsD2Zrk-eqEUR94UtYotqJ"
$ws.Cells.Item(5, 7).Value = "Mechanical Engineering "
$ws.Cells.Item(5, 8).Value = "Kangra"
$ws.Cells.Item(5, 9).Value = "Himachal Pradesh "
$ws.Cells.Item(5, 10).Value = "Batch 2025"

$ws.Cells.Item(6, 1).Value = "Dharuva Thakur "
$ws.Cells.Item(6, 2).Value = "first"
$ws.Cells.Item(6, 3).Value = "Volunteer Member"
$ws.Cells.Item(6, 4).Value = "https://www.linkedin.com/in/dharuva-thakur-83576122a"
$ws.Cells.Item(6, 5).Value = "../members/dharuva.webp"
$ws.Cells.Item(6, 6).Value = "18hlg1lzF81JB-DJKg36G5SHzOM2zB7Xz"
$ws.Cells.Item(6, 7).Value = "Mathematics and Scientific Computing"
$ws.Cells.Item(6, 8).Value = "Mandi"
$ws.Cells.Item(6, 9).Value = "Himachal Pradesh "
$ws.Cells.Item(6, 10).Value = "Batch 2025"

$ws.Cells.Item(7, 1).Value = "Ekansh Verma"
$ws.Cells.Item(7, 2).Value = "first"
$ws.Cells.Item(7, 3).Value = "Volunteer Member"
$ws.Cells.Item(7, 4).Value = "https://www.linkedin.com/in/ekansh-verma-469557228"
$ws.Cells.Item(7, 5).Value = "../members/ekansh.webp"
$ws.Cells.Item(7, 6).Value = "1upS_dc_SLM0CsaA__dVnXqApgS4eFKD1"
$ws.Cells.Item(7, 7).Value = "Electrical Engineering"
$ws.Cells.Item(7, 8).Value = "Kullu"
$ws.Cells.Item(7, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(7, 10).Value = "First Year UG at NiTH"

$ws.Cells.Item(8, 1).Value = "Gargi Dhawan"
$ws.Cells.Item(8, 2).Value = "first"
$ws.Cells.Item(8, 3).Value = "Volunteer Member"
$ws.Cells.Item(8, 4).Value = "https://www.linkedin.com/in/gargi-dhawan-7706a722a"
$ws.Cells.Item(8, 5).Value = "../members/gargi.webp"
$ws.Cells.Item(8, 6).Value = "1HD-obHJOhcS0BYxzv4SRHdRGRnCJ4d-H"
$ws.Cells.Item(8, 7).Value = "Electronics and Communication Engineering"
$ws.Cells.Item(8, 8).Value = "Shimla"
$ws.Cells.Item(8, 9).Value = "Himacha pradesh"
$ws.Cells.Item(8, 10).Value = "Batch 2025"

$ws.Cells.Item(9, 1).Value = "Hardik Sachdeva"
$ws.Cells.Item(9, 2).Value = "first"
$ws.Cells.Item(9, 3).Value = "Volunteer Member"
$ws.Cells.Item(9, 4).Value = "https://www.linkedin.com/in/hardik-sachdeva-a69987217"
$ws.Cells.Item(9, 5).Value = "../members/hardik.webp"
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(9, 7).Value = "Mathematics and Scientific Computing"
$ws.Cells.Item(9, 8).Value = "Nainital"
$ws.Cells.Item(9, 9).Value = "Uttarakhand"
$ws.Cells.Item(9, 10).Value = "Pursuing BTech in Mathematics And Computing (2025), Interested in Tech. And Financial Stuff. Travel Freak"

$ws.Cells.Item(10, 1).Value = "Jeevak Sangodkar"
$ws.Cells.Item(10, 2).Value = "first"
$ws.Cells.Item(10, 3).Value = "Volunteer Member"
$ws.Cells.Item(10, 4).Value = "https://www.linkedin.com/in/jeevak-sangodkar-919653228/"
$ws.Cells.Item(10, 5).Value = "../members/jeevak.webp"
$ws.Cells.Item(10, 6).Value = "16_rAwwruLc2BS3XuKMGjF2K70foAL39d"
$ws.Cells.Item(10, 7).Value = "Computer science and Engineering"
$ws.Cells.Item(10, 8).Value = "Nagpur"
$ws.Cells.Item(10, 9).Value = "Maharashtra"
$ws.Cells.Item(10, 10).Value = "Batch 2025"

$ws.Cells.Item(11, 1).Value = "Kanika Sharma "
$ws.Cells.Item(11, 2).Value = "first"
$ws.Cells.Item(11, 3).Value = "Volunteer  Member"
$ws.Cells.Item(11, 4).Value = "https://www.linkedin.com/in/kanika-sharma-aa274b22a"
$ws.Cells.Item(11, 5).Value = "../members/kanika.webp"
$ws.Cells.Item(11, 6).Value = "1cMSLUYaTd_-_d0bWjLBF-VVSiQcmPDaX"
$ws.Cells.Item(11, 7).Value = "Electrical Engineering"
$ws.Cells.Item(11, 8).Value = "Kangra"
$ws.Cells.Item(11, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(11, 10).Value = "First year undergrad at NITH; technology enthusiast"

$ws.Cells.Item(12, 1).Value = "Kuldeep Kumar"
$ws.Cells.Item(12, 2).Value = "first"
$ws.Cells.Item(12, 3).Value = "Volunteer Member"
$ws.Cells.Item(12, 5).Value = "../members/kuldeep.webp"
$ws.Cells.Item(12, 6).Value = "156x7nF6Vjqlqz0aV1u_krIV2F0B-lvpm"
$ws.Cells.Item(12, 7).Value = "Computer science and Engineering"
$ws.Cells.Item(12, 8).Value = "Solan"
$ws.Cells.Item(12, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(12, 10).Value = "Batch 2025"

$ws.Cells.Item(13, 1).Value = "Madhukesh Singh"
$ws.Cells.Item(13, 2).Value = "first"
$ws.Cells.Item(13, 3).Value = "Volunteer Member"
$ws.Cells.Item(13, 4).Value = "https://www.linkedin.com/in/madhukesh-singh-195618233/"
$ws.Cells.Item(13, 5).Value = "../members/madhukesh.webp"
$ws.Cells.Item(13, 6).Value = "1G1H5v8mdSJ_Us4i1IQtIrSrkYvOBWHWm"
$ws.Cells.Item(13, 7).Value = "Chemical Engineering"
$ws.Cells.Item(13, 8).Value = "Patna"
$ws.Cells.Item(13, 9).Value = "Bihar"
$ws.Cells.Item(13, 10).Value = "Batch 2025"

$ws.Cells.Item(14, 1).Value = "Manik Singh "
$ws.Cells.Item(14, 2).Value = "first"
$ws.Cells.Item(14, 3).Value = "Volunteer Member"
$ws.Cells.Item(14, 5).Value = "../members/manik.webp"
$ws.Cells.Item(14, 6).Value = "1ijjGG-_dZMlx93qwdHh_MsDrI7vzqgGe"
$ws.Cells.Item(14, 7).Value = "Electronics and Communication Engineering"
$ws.Cells.Item(14, 8).Value = "Kangra"
$ws.Cells.Item(14, 9).Value = "Himachal Pradesh "
$ws.Cells.Item(14, 10).Value = "Batch 2025"

$ws.Cells.Item(15, 1).Value = "Mehul Aggarwal"
$ws.Cells.Item(15, 2).Value = "first"
$ws.Cells.Item(15, 3).Value = "Volunteer Member"
$ws.Cells.Item(15, 4).Value = "https://www.linkedin.com/in/mehul-aggarwal-47285421b/"
$ws.Cells.Item(15, 5).Value = "../members/mehul.webp"
$ws.Cells.Item(15, 6).Value = "1_J6_Pqr7pITVZdau7N3K8iZZXftah1kj"
$ws.Cells.Item(15, 7).Value = "Mathematics and Scientific Computing"
$ws.Cells.Item(15, 8).Value = "Haridwar"
$ws.Cells.Item(15, 9).Value = "Uttarakhand"
$ws.Cells.Item(15, 10).Value = "A student at NIT Hamirpur currently pursuing Mathematics and Computing; learning ,trying and testing new technologies and skills everyday."

$ws.Cells.Item(16, 1).Value = "Navdeep Kaur"
$ws.Cells.Item(16, 2).Value = "first"
$ws.Cells.Item(16, 3).Value = "Volunteer Member"
$ws.Cells.Item(16, 4).Value = "https://www.linkedin.com/in/navdeep-kaur-44375022a/"
$ws.Cells.Item(16, 5).Value = "../members/navdeep.webp"
$ws.Cells.Item(16, 6).Value = "1xKhiNxAbeV8uBDzZRxc6_rxo-opXstOS"
$ws.Cells.Item(16, 7).Value = "Computer science and Engineering"
$ws.Cells.Item(16, 8).Value = "Bilaspur"
$ws.Cells.Item(16, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(16, 10).Value = "Batch 2025"

$ws.Cells.Item(17, 1).Value = "Shashank Shekhar"
$ws.Cells.Item(17, 2).Value = "first"
$ws.Cells.Item(17, 3).Value = "Volenteer Member"
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = "../members/shashank.webp"
$ws.Cells.Item(17, 6).Value = ""
$ws.Cells.Item(17, 7).Value = ""
$ws.Cells.Item(17, 8).Value = "Jamui"
$ws.Cells.Item(17, 9).Value = "Bihar"
$ws.Cells.Item(17, 10).Value = "Batch 2025"

$ws.Cells.Item(18, 1).Value = "Shariq Verma"
$ws.Cells.Item(18, 2).Value = "first"
$ws.Cells.Item(18, 3).Value = "Volunteer  Member"
$ws.Cells.Item(18, 4).Value = "https://www.linkedin.com/in/shariq-verma-94a75122a"
$ws.Cells.Item(18, 5).Value = "../members/Shariq.webp"
$ws.Cells.Item(18, 6).Value = "1KC9wcQ5DgKU31KofULgTIa5VIsbtrttt"
$ws.Cells.Item(18, 7).Value = "Civil Engineering"
$ws.Cells.Item(18, 8).Value = "Shimla"
$ws.Cells.Item(18, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(18, 10).Value = "A Civil Engineering Undergraduate interested in core and currently working on core based technical skills. "

$ws.Cells.Item(19, 1).Value = "Siya Rana"
$ws.Cells.Item(19, 2).Value = "first"
$ws.Cells.Item(19, 3).Value = "Volunteer Member"
$ws.Cells.Item(19, 4).Value = "https://www.linkedin.com/in/siya-rana-b07639221"
$ws.Cells.Item(19, 5).Value = "../members/siya.webp"
$ws.Cells.Item(19, 6).Value = "1K4LpZE6yPe1EAsxK6zL7RwGDRhvTk97t"
$ws.Cells.Item(19, 7).Value = "Electronics and Communication Engineering"
$ws.Cells.Item(19, 8).Value = "Kangra"
$ws.Cells.Item(19, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(19, 10).Value = "First year student at Nit Hamirpur, pursuing Electronics and Communication Engineering."

$ws.Cells.Item(20, 1).Value = "Swastik Sharma"
$ws.Cells.Item(20, 2).Value = "first"
$ws.Cells.Item(20, 3).Value = "Volunteer Member"
$ws.Cells.Item(20, 4).Value = "https://www.linkedin.com/in/swastkk"
$ws.Cells.Item(20, 5).Value = "../members/Swastik.webp"
$ws.Cells.Item(20, 6).Value = "15p06FU4eSbaUxpU2A1mq4bD9BS5PxPxu"
$ws.Cells.Item(20, 7).Value = "Civil Engineering"
$ws.Cells.Item(20, 8).Value = "UNA"
$ws.Cells.Item(20, 9).Value = "Himachal Pradesh"
$ws.Cells.Item(20, 10).Value = "A CE Undergrad at NIT Hamirpur currently learning new Web Technologies at Top pace."
